# Assignments.xlsx - add "Assignment_6" row and restyle the link cells
# (blue, non-underlined custom color instead of the default themed
# Hyperlink look) to match the rest of the "link" column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new "Assignment_6" row -----------------------------------
# Clone row 5's look (full thin border, vertical-center alignment) into
# row 7, then overwrite the three values for the new assignment.
$ws.Range("A5:C5").Copy() | Out-Null
$ws.Range("A7:C7").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Application.CutCopyMode = $false

$ws.Rows.Item(7).RowHeight = 31.5

$ws.Range("A7").Value = "Assignment_6"
$ws.Range("B7").Value = "https://github.com/Vasanth30e/Assignments_Phase2/tree/master/Assignment_6"
$ws.Range("C7").Value = 45156

# --- Re-color the GitHub-link cells ------------------------------------
# Drop the underline and switch from the themed Hyperlink blue to a
# custom blue (RGB 0,112,192 -> 0x0070C0), applied to each link cell
# (including the freshly added one).
$linkColor = 12611584   # RGB(0, 112, 192)

$ws.Range("B3").Font.Underline = $false
$ws.Range("B3").Font.Color = $linkColor

$ws.Range("B4").Font.Underline = $false
$ws.Range("B4").Font.Color = $linkColor

$ws.Range("B5").Font.Underline = $false
$ws.Range("B5").Font.Color = $linkColor

$ws.Range("B6").Font.Underline = $false
$ws.Range("B6").Font.Color = $linkColor

$ws.Range("B7").Font.Underline = $false
$ws.Range("B7").Font.Color = $linkColor

# --- Match the saved cursor position from the source workbook ---------
$ws.Range("B14").Select() | Out-Null
